$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2023976666666667
$ws.Range("N2").Value = 0.6071930000000001
$ws.Range("O2").Value = 0.03663970451354832
$ws.Range("P2").Value = 0.03663970451354832
$ws.Range("Q2").Value = 0.06503313640144445
$ws.Range("R2").Value = 0.5852982276130001
$ws.Range("S2").Value = 0.02729823683105628
$ws.Range("T2").Value = 0.02729823683105628

# Row 3 (FAPs -> FAPs)
$ws.Range("O3").Value = 0.08641717548188978
$ws.Range("P3").Value = 0.08641717548188979
$ws.Range("S3").Value = 0.06438470380412792
$ws.Range("T3").Value = 0.06438470380412793

# Row 4 (FAPs -> MuSCs)
$ws.Range("O4").Value = 0.876943120004562
$ws.Range("P4").Value = 0.876943120004562
$ws.Range("S4").Value = 0.6533622826679178
$ws.Range("T4").Value = 0.6533622826679178

# Row 5 (MuSCs -> ECs)
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2023976666666667
$ws.Range("N5").Value = 0.6071930000000001
$ws.Range("O5").Value = 0.03663970451354832
$ws.Range("P5").Value = 0.03663970451354832
$ws.Range("Q5").Value = 0.02225436557477778
$ws.Range("R5").Value = 0.200289290173
$ws.Range("S5").Value = 0.009341467682492036
$ws.Range("T5").Value = 0.009341467682492036

# Row 6 (MuSCs -> FAPs)
$ws.Range("O6").Value = 0.08641717548188978
$ws.Range("P6").Value = 0.08641717548188979
$ws.Range("S6").Value = 0.02203247167776185
$ws.Range("T6").Value = 0.02203247167776186

# Row 7 (MuSCs -> MuSCs)
$ws.Range("O7").Value = 0.876943120004562
$ws.Range("P7").Value = 0.876943120004562
$ws.Range("S7").Value = 0.2235808373366441
$ws.Range("T7").Value = 0.2235808373366441
